# Fruta / hortaliza, semanal
# Insert a new weekly record at row 66, pushing the existing rows 66:109
# down to 67:110, and populate the new row with the latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("66:66").Insert()

$ws.Range("A66").Value2 = 6
$ws.Range("B66").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C66").Value2 = "Metropolitana"
$ws.Range("D66").Value2 = 44488
$ws.Range("E66").Value2 = 13
$ws.Range("F66").Value2 = 100112001
$ws.Range("G66").Value2 = "Berenjena"
$ws.Range("H66").Value2 = "Sin especificar"
$ws.Range("I66").Value2 = "Primera"
$ws.Range("J66").Value2 = 210
$ws.Range("K66").Value2 = 7000
$ws.Range("L66").Value2 = 8000
$ws.Range("M66").Value2 = 7619
$ws.Range("N66").Value2 = "`$/caja 60 unidades"
$ws.Range("O66").Value2 = "Provincia de Huasco"
$ws.Range("P66").Value2 = 127
$ws.Range("Q66").Value2 = 60
$ws.Range("R66").Value2 = "Hortaliza"
